$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted numeric-looking strings.
# Force text number format so Excel does not auto-convert these values to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '39.219.76'
$ws.Range("E2").Value = '  -1.91%  '

$ws.Range("D3").Value = '2.199.60'
$ws.Range("E3").Value = '  -5.51%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '294.93'
$ws.Range("E5").Value = '  -3.89%  '

$ws.Range("D6").Value = '81.07'
$ws.Range("E6").Value = '  -4.20%  '

$ws.Range("E7").Value = '  -3.66%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = '0.467'
$ws.Range("E9").Value = '  -3.07%  '

$ws.Range("D10").Value = '0.0768'
$ws.Range("E10").Value = '  -5.58%  '

$ws.Range("D11").Value = '29.05'
$ws.Range("E11").Value = '  -3.07%  '

$ws.Range("D12").Value = '46.76'
$ws.Range("E12").Value = '  -11.15%  '

$ws.Range("E13").Value = '  -2.50%  '

$ws.Range("D14").Value = '2.529.10'
$ws.Range("E14").Value = '  -5.92%  '

$ws.Range("D15").Value = '6.22'
$ws.Range("E15").Value = '  -2.72%  '

$ws.Range("D16").Value = '13.92'
$ws.Range("E16").Value = '  -4.95%  '

$ws.Range("D17").Value = '2.199.06'
$ws.Range("E17").Value = '  -5.09%  '

$ws.Range("D18").Value = '0.709'
$ws.Range("E18").Value = '  -5.24%  '

$ws.Range("D19").Value = '39.120.00'
$ws.Range("E19").Value = '  -2.15%  '

$ws.Range("D20").Value = '0.0₃0869'
$ws.Range("E20").Value = '  -3.43%  '

$ws.Range("D21").Value = '5.70'
$ws.Range("E21").Value = '  -6.03%  '

$ws.Range("D22").Value = '64.55'
$ws.Range("E22").Value = '  -4.36%  '

$ws.Range("D23").Value = '10.24'
$ws.Range("E23").Value = '  -3.42%  '

$ws.Range("D24").Value = '225.84'
$ws.Range("E24").Value = '  -3.90%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  -6.26%  '

$ws.Range("D27").Value = '1.80'
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").Value = '22.49'
$ws.Range("E28").Value = '  -3.51%  '

$ws.Range("E29").Value = '  -1.55%  '

$ws.Range("D30").Value = '9.02'
$ws.Range("E30").Value = '  -2.23%  '

$ws.Range("D31").Value = '149.32'
$ws.Range("E31").Value = '  -1.66%  '

$ws.Range("D32").Value = '31.53'
$ws.Range("E32").Value = '  -10.44%  '

$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.31%  '

$ws.Range("D34").Value = '4.79'
$ws.Range("E34").Value = '  -5.97%  '

$ws.Range("E35").Value = '  -4.43%  '

$ws.Range("D36").Value = '0.0693'
$ws.Range("E36").Value = '  -3.88%  '

$ws.Range("D37").Value = '0.110'
$ws.Range("E37").Value = '  -3.25%  '

$ws.Range("D38").Value = '15.29'
$ws.Range("E38").Value = '  -2.58%  '

$ws.Range("D39").Value = '0.0958'
$ws.Range("E39").Value = '  -3.81%  '

$ws.Range("D40").Value = '2.62'
$ws.Range("E40").Value = '  -4.35%  '

$ws.Range("D41").Value = '1.65'
$ws.Range("E41").Value = '  -3.03%  '

$ws.Range("E42").Value = '  -5.39%  '

$ws.Range("D43").Value = '1.900.07'
$ws.Range("E43").Value = '  -1.87%  '

$ws.Range("D44").Value = '2.04'
$ws.Range("E44").Value = '  -9.93%  '

$ws.Range("E45").Value = '  -2.47%  '

$ws.Range("D46").Value = '8.99'
$ws.Range("E46").Value = '  -2.80%  '

$ws.Range("D47").Value = '15.99'
$ws.Range("E47").Value = '  -8.96%  '

$ws.Range("D48").Value = '2.60'
$ws.Range("E48").Value = '  -2.57%  '

$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").Value = '71.25'
$ws.Range("E49").Value = '  +0.54%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.405.67'
$ws.Range("E50").Value = '  -5.97%  '

$ws.Range("D51").Value = '87.17'
$ws.Range("E51").Value = '  -5.91%  '
